$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.637.22'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').Value = '2.526.11'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''316.18'
$ws.Range('E5').Value = '  +4.04%  '
$ws.Range('D6').Value = '''94.83'
$ws.Range('E6').Value = '  -3.43%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -1.33%  '
$ws.Range('D10').Value = '''36.29'
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('D11').Value = '''0.0811'
$ws.Range('E11').Value = '  -2.01%  '
$ws.Range('D12').Value = '''7.69'
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').Value = '''0.114'
$ws.Range('E13').Value = '  -2.54%  '
$ws.Range('D14').Value = '2.906.40'
$ws.Range('E14').Value = '  -1.30%  '
$ws.Range('D15').Value = '''15.53'
$ws.Range('E15').Value = '  +2.92%  '
$ws.Range('D16').Value = '2.503.43'
$ws.Range('E16').Value = '  -3.03%  '
$ws.Range('D17').Value = '''0.863'
$ws.Range('E17').Value = '  -1.74%  '
$ws.Range('D18').Value = '42.714.42'
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('E19').Value = '  -5.86%  '
$ws.Range('D20').Value = '0.0₃0967'
$ws.Range('E20').Value = '  -2.98%  '
$ws.Range('D21').Value = '''6.54'
$ws.Range('E21').Value = '  -0.88%  '
$ws.Range('D22').Value = '''71.24'
$ws.Range('E22').Value = '  -1.06%  '
$ws.Range('D23').Value = '''254.28'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').Value = '''2.98'
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('E25').Value = '  -2.32%  '
$ws.Range('D26').Value = '''27.61'
$ws.Range('E26').Value = '  -1.59%  '
$ws.Range('D27').Value = '''0.992'
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('D28').Value = '''2.32'
$ws.Range('E28').Value = '  +10.98%  '
$ws.Range('D29').Value = '''39.29'
$ws.Range('E29').Value = '  +3.76%  '
$ws.Range('D30').Value = '''10.05'
$ws.Range('E30').Value = '  -1.78%  '
$ws.Range('D31').Value = '''5.90'
$ws.Range('E31').Value = '  -4.29%  '
$ws.Range('D32').Value = '''156.12'
$ws.Range('E32').Value = '  -1.65%  '
$ws.Range('D33').Value = '''19.91'
$ws.Range('E33').Value = '  +1.41%  '
$ws.Range('D34').Value = '''3.33'
$ws.Range('E34').Value = '  +0.66%  '
$ws.Range('D35').Value = '''2.09'
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('D36').Value = '''0.0784'
$ws.Range('E36').Value = '  -2.43%  '
$ws.Range('E37').Value = '  -1.24%  '
$ws.Range('E38').Value = '  -3.19%  '
$ws.Range('D39').Value = '''24.73'
$ws.Range('E39').Value = '  -3.16%  '
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('D41').Value = '''2.17'
$ws.Range('E41').Value = '  +3.22%  '
$ws.Range('D42').Value = '''3.85'
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('D43').Value = '''3.37'
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '''0.0303'
$ws.Range('E44').Value = '  -1.11%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '''1.00'
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('D46').Value = '2.053.37'
$ws.Range('E46').Value = '  -1.72%  '
$ws.Range('D47').Value = '''86.36'
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').Value = '''8.83'
$ws.Range('E48').Value = '  -1.72%  '
$ws.Range('D49').Value = '2.760.19'
$ws.Range('E49').Value = '  -1.53%  '
$ws.Range('D50').Value = '''74.23'
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('E51').Value = '  -0.53%  '
